$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 90999.5
$ws.Range("J7").Value = 90999
$ws.Range("L7").Value = 90999
$ws.Range("N7").Value = -91223
$ws.Range("H14").Value = 90999.5
$ws.Range("J14").Value = 90999
$ws.Range("L14").Value = 90999
$ws.Range("N14").Value = -91381
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = ""
$ws.Range("H111").Value = 3591.1667
$ws.Range("I111").Value = 2428.875
$ws.Range("K111").Value = 7286.625
$ws.Range("M111").Value = -4219.625
$ws.Range("H113").Value = 3509.75
$ws.Range("I113").Value = 2349
$ws.Range("K113").Value = 2349
$ws.Range("M113").Value = 905
$ws.Range("H116").Value = 8750
$ws.Range("I116").Value = 7500
$ws.Range("K116").Value = 7500
$ws.Range("M116").Value = -4058
$ws.Range("H138").Value = 6979.5386
$ws.Range("I138").Value = 11167.4
$ws.Range("K138").Value = 33502.2
$ws.Range("M138").Value = -28362.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1644.7307
$ws.Range("I61").Value = 1457.6364
$ws.Range("K61").Value = 1457.6364
$ws.Range("M61").Value = -1245.6364
$ws.Range("H74").Value = 1262.7435
$ws.Range("I74").Value = 831.6857
$ws.Range("K74").Value = 831.6857
$ws.Range("M74").Value = 42.3143
$ws.Range("H77").Value = 1262.7435
$ws.Range("I77").Value = 831.6857
$ws.Range("K77").Value = 4158.4285
$ws.Range("M77").Value = 209.5715
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H132").Value = 1293.3793
$ws.Range("I132").Value = 1293.3793
$ws.Range("K132").Value = 3880.1379
$ws.Range("M132").Value = -1350.1379
$ws.Range("H136").Value = 1644.7307
$ws.Range("I136").Value = 1457.6364
$ws.Range("K136").Value = 4372.9092
$ws.Range("M136").Value = -1822.9092

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1557.5
$ws.Range("I86").Value = 1769.5
$ws.Range("J86").Value = 497.5
$ws.Range("K86").Value = 1769.5
$ws.Range("L86").Value = 497.5
$ws.Range("M86").Value = -646.5
$ws.Range("N86").Value = -2743.5
$ws.Range("H89").Value = 1557.5
$ws.Range("I89").Value = 1769.5
$ws.Range("J89").Value = 497.5
$ws.Range("K89").Value = 8847.5
$ws.Range("L89").Value = 2487.5
$ws.Range("M89").Value = -3231.5
$ws.Range("N89").Value = -13719.5
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1650
$ws.Range("K107").Value = 1650
$ws.Range("M107").Value = 270
$ws.Range("H134").Value = 2993.7144
$ws.Range("I134").Value = 2774.4375
$ws.Range("K134").Value = 8323.3125
$ws.Range("M134").Value = -5788.3125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1401
$ws.Range("I16").Value = 1005
$ws.Range("K16").Value = 1005
$ws.Range("M16").Value = -718
$ws.Range("H31").Value = 4182.037
$ws.Range("I31").Value = 2582.6296
$ws.Range("J31").Value = 5781.4443
$ws.Range("K31").Value = 2582.6296
$ws.Range("L31").Value = 5781.4443
$ws.Range("M31").Value = -2287.6296
$ws.Range("N31").Value = -6371.4443
$ws.Range("H34").Value = 4182.037
$ws.Range("I34").Value = 2582.6296
$ws.Range("J34").Value = 5781.4443
$ws.Range("K34").Value = 2582.6296
$ws.Range("L34").Value = 5781.4443
$ws.Range("M34").Value = -2380.6296
$ws.Range("N34").Value = -6185.4443
$ws.Range("H51").Value = 26653.8
$ws.Range("I51").Value = 20817.5
$ws.Range("J51").Value = 49999
$ws.Range("K51").Value = 20817.5
$ws.Range("L51").Value = 49999
$ws.Range("M51").Value = -20081.5
$ws.Range("N51").Value = -51471.5
$ws.Range("H61").Value = 26653.8
$ws.Range("I61").Value = 20817.5
$ws.Range("J61").Value = 49999
$ws.Range("K61").Value = 20817.5
$ws.Range("L61").Value = 49999
$ws.Range("M61").Value = -20469.5
$ws.Range("N61").Value = -50695
$ws.Range("H113").Value = 1401
$ws.Range("I113").Value = 1005
$ws.Range("K113").Value = 1005
$ws.Range("M113").Value = 1165
$ws.Range("H132").Value = 2069.577
$ws.Range("I132").Value = 1529.3617
$ws.Range("K132").Value = 4588.0851
$ws.Range("M132").Value = -2058.0851
$ws.Range("H134").Value = 2902.25
$ws.Range("I134").Value = 2282.7896
$ws.Range("K134").Value = 6848.3688
$ws.Range("M134").Value = -4313.3688

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5800
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 2433.7778
$ws.Range("I132").Value = 1873.8889
$ws.Range("J132").Value = 3553.5557
$ws.Range("K132").Value = 5621.6667
$ws.Range("L132").Value = 10660.6671
$ws.Range("M132").Value = -3091.6667
$ws.Range("N132").Value = -15720.6671

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15066.333
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -605
$ws.Range("H27").Value = 15066.333
$ws.Range("I27").Value = 900
$ws.Range("K27").Value = 900
$ws.Range("M27").Value = -793
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H136").Value = 3041.2856
$ws.Range("I136").Value = 849.25
$ws.Range("K136").Value = 2547.75
$ws.Range("M136").Value = 2.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H132").Value = 1865.6666
$ws.Range("J132").Value = 1932
$ws.Range("L132").Value = 5796
$ws.Range("N132").Value = -10856
$ws.Range("H136").Value = 1146.1
$ws.Range("I136").Value = 1151.8148
$ws.Range("K136").Value = 3455.4444
$ws.Range("M136").Value = -905.4444000000003
